# Rename worksheet tab to reflect new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-11-07"

# Row 12 (October) - T,U,V columns (2021 year block)
$ws.Range("T12").Value = 3
$ws.Range("U12").Value = 193
$ws.Range("V12").Value = 0.0153

# Row 13 (November through ...) - update label and values
$ws.Range("A13").Value = "November (through 11-07)"
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 17
$ws.Range("I13").Value = 23
$ws.Range("J13").Value = 0.0417
$ws.Range("N13").Value = 1
$ws.Range("P13").Value = 0.1
$ws.Range("R13").Value = 44
$ws.Range("S13").Value = 0.0222
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 47
$ws.Range("V13").Value = 0.0208

# Row 14 (Total) - updated aggregate values
$ws.Range("C14").Value = 234
$ws.Range("D14").Value = 0.1203
$ws.Range("F14").Value = 451
$ws.Range("G14").Value = 0.1034
$ws.Range("I14").Value = 672
$ws.Range("J14").Value = 0.08450000000000001
$ws.Range("N14").Value = 49
$ws.Range("P14").Value = 0.09959999999999999
$ws.Range("R14").Value = 1047
$ws.Range("S14").Value = 0.0499
$ws.Range("T14").Value = 85
$ws.Range("U14").Value = 1407
$ws.Range("V14").Value = 0.057
